$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.168.16'
$ws.Range("E2").Value = '  +2.80%  '

$ws.Range("D3").Value = '2.266.92'
$ws.Range("E3").Value = '  +2.49%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = "'270.53"
$ws.Range("E5").Value = '  +5.31%  '

$ws.Range("D6").Value = "'87.60"
$ws.Range("E6").Value = '  +13.49%  '

$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = '  +1.26%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = "'0.614"
$ws.Range("E9").Value = '  +2.93%  '

$ws.Range("D10").Value = "'45.79"
$ws.Range("E10").Value = '  +6.86%  '

$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = '  +2.52%  '

$ws.Range("D12").Value = "'7.65"
$ws.Range("E12").Value = '  +8.74%  '

$ws.Range("E13").Value = '  +2.43%  '

$ws.Range("D14").Value = '2.607.99'
$ws.Range("E14").Value = '  +2.69%  '

$ws.Range("D15").Value = "'15.09"
$ws.Range("E15").Value = '  +4.38%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.258.57'
$ws.Range("E16").Value = '  +2.26%  '

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = "'0.807"
$ws.Range("E17").Value = '  +2.91%  '

$ws.Range("D18").Value = '44.064.61'
$ws.Range("E18").Value = '  +2.77%  '

$ws.Range("D19").Value = "'0.0000104"
$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("D20").Value = "'6.08"
$ws.Range("E20").Value = '  +1.84%  '

$ws.Range("D21").Value = "'70.66"
$ws.Range("E21").Value = '  -0.70%  '

$ws.Range("D22").Value = "'2.40"
$ws.Range("E22").Value = '  +2.43%  '

$ws.Range("D23").Value = "'235.43"
$ws.Range("E23").Value = '  +2.28%  '

$ws.Range("D24").Value = "'8.97"
$ws.Range("E24").Value = '  -2.77%  '

$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = "'2.53"
$ws.Range("E26").Value = '  +15.22%  '

$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = '  +1.78%  '

$ws.Range("E28").Value = '  +6.69%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = "'2.33"
$ws.Range("E29").Value = '  +5.46%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = "'39.75"
$ws.Range("E30").Value = '  -6.33%  '

$ws.Range("D31").Value = "'174.71"
$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'21.06"
$ws.Range("E32").Value = '  +3.43%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.0904"
$ws.Range("E33").Value = '  +4.08%  '

$ws.Range("D34").Value = "'5.39"
$ws.Range("E34").Value = '  +2.97%  '

$ws.Range("E35").Value = '  +1.48%  '

$ws.Range("E36").Value = '  +5.97%  '

$ws.Range("D37").Value = "'0.0353"
$ws.Range("E37").Value = '  -2.70%  '

$ws.Range("D38").Value = "'4.40"
$ws.Range("E38").Value = '  +0.44%  '

$ws.Range("D39").Value = "'3.40"
$ws.Range("E39").Value = '  +16.35%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = "'2.24"
$ws.Range("E40").Value = '  +5.98%  '

$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").Value = "'12.82"
$ws.Range("E41").Value = '  -1.71%  '

$ws.Range("D42").Value = "'64.72"
$ws.Range("E42").Value = '  +5.73%  '

$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").Value = "'5.49"
$ws.Range("E43").Value = '  +4.07%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.205"
$ws.Range("E44").Value = '  +1.14%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'102.17"
$ws.Range("E45").Value = '  -0.82%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = "'0.0997"
$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'8.44"
$ws.Range("E47").Value = '  -0.35%  '

$ws.Range("E48").Value = '  +7.80%  '

$ws.Range("E49").Value = '  +3.02%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = "'1.54"
$ws.Range("E50").Value = '  +5.00%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.437"
$ws.Range("E51").Value = '  -6.44%  '
